# Update "想去人数" (column F) values on the "展览" and "全部类型" worksheets.
$wb = $excel.ActiveWorkbook

$updates = @{
    4  = 1557
    5  = 595
    6  = 1086
    7  = 11286
    8  = 12
    10 = 440
    11 = 337
    12 = 1082
    14 = 12294
    15 = 12930
    22 = 80
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
